# Auto-generated edit script: updates numeric cell values in the
# ALC / ARM / CRP / CUL / GSM / LTW / WVR sheets per the scheduled-runner diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 945.1786
$ws.Range("I98").Value = 945.1786
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 945.1786
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 552.8214
$ws.Range("N98").ClearContents()
$ws.Range("H100").Value = 5954007.5
$ws.Range("J100").Value = 1787.4117
$ws.Range("L100").Value = 1787.4117
$ws.Range("N100").Value = -2869.4117
$ws.Range("H101").Value = 15160
$ws.Range("I101").Value = 250
$ws.Range("J101").Value = 25100
$ws.Range("K101").Value = 750
$ws.Range("L101").Value = 75300
$ws.Range("M101").Value = 872
$ws.Range("N101").Value = -78544
$ws.Range("H105").Value = 31500
$ws.Range("J105").Value = 31500
$ws.Range("L105").Value = 31500
$ws.Range("N105").Value = -38488
$ws.Range("H106").Value = 190477540
$ws.Range("I106").Value = 55557136
$ws.Range("K106").Value = 55557136
$ws.Range("M106").Value = -55556505
$ws.Range("H107").Value = 41676068
$ws.Range("I107").Value = 50001280
$ws.Range("J107").Value = 50000
$ws.Range("K107").Value = 50001280
$ws.Range("L107").Value = 50000
$ws.Range("M107").Value = -49999360
$ws.Range("N107").Value = -53840
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H111").Value = 168516.83
$ws.Range("I111").Value = 1500
$ws.Range("J111").Value = 252025.25
$ws.Range("K111").Value = 4500
$ws.Range("L111").Value = 756075.75
$ws.Range("M111").Value = -1433
$ws.Range("N111").Value = -762209.75
$ws.Range("H112").Value = 6850390
$ws.Range("J112").Value = 7354044.5
$ws.Range("L112").Value = 22062133.5
$ws.Range("N112").Value = -22064349.5
$ws.Range("H113").Value = 8700
$ws.Range("I113").Value = 3050
$ws.Range("J113").Value = 20000
$ws.Range("K113").Value = 3050
$ws.Range("L113").Value = 20000
$ws.Range("M113").Value = 204
$ws.Range("N113").Value = -26508
$ws.Range("H115").Value = 698.5
$ws.Range("I115").Value = 698.5
$ws.Range("K115").Value = 2095.5
$ws.Range("M115").Value = -528.5
$ws.Range("H116").Value = 8162.5
$ws.Range("I116").Value = 10150
$ws.Range("J116").Value = 2200
$ws.Range("K116").Value = 10150
$ws.Range("L116").Value = 2200
$ws.Range("M116").Value = -6708
$ws.Range("N116").Value = -9084
$ws.Range("H117").Value = 27800
$ws.Range("J117").Value = 27800
$ws.Range("L117").Value = 27800
$ws.Range("N117").Value = -36978
$ws.Range("H118").Value = 657
$ws.Range("I118").Value = 566.5
$ws.Range("J118").Value = 1200
$ws.Range("K118").Value = 1699.5
$ws.Range("L118").Value = 3600
$ws.Range("M118").Value = -42.5
$ws.Range("N118").Value = -6914
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H122").Value = 945.1786
$ws.Range("I122").Value = 945.1786
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2835.5358
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -385.5357999999997
$ws.Range("N122").ClearContents()
$ws.Range("H123").Value = 35749.5
$ws.Range("J123").Value = 35749.5
$ws.Range("L123").Value = 35749.5
$ws.Range("N123").Value = -45549.5
$ws.Range("H124").Value = 29000
$ws.Range("J124").Value = 29000
$ws.Range("L124").Value = 29000
$ws.Range("N124").Value = -38820
$ws.Range("H125").Value = 6227.091
$ws.Range("J125").Value = 6699.8
$ws.Range("L125").Value = 60298.2
$ws.Range("N125").Value = -65218.2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 10446.077
$ws.Range("I61").Value = 12689.9
$ws.Range("J61").Value = 2966.6667
$ws.Range("K61").Value = 12689.9
$ws.Range("L61").Value = 2966.6667
$ws.Range("M61").Value = -12477.9
$ws.Range("N61").Value = -3390.6667
$ws.Range("H63").Value = 125002760
$ws.Range("I63").Value = 142859580
$ws.Range("K63").Value = 142859580
$ws.Range("M63").Value = -142858894
$ws.Range("H66").Value = 125002760
$ws.Range("I66").Value = 142859580
$ws.Range("K66").Value = 714297900
$ws.Range("M66").Value = -714294468
$ws.Range("H74").Value = 1650.9242
$ws.Range("I74").Value = 1560.2322
$ws.Range("J74").Value = 2158.8
$ws.Range("K74").Value = 1560.2322
$ws.Range("L74").Value = 2158.8
$ws.Range("M74").Value = -686.2321999999999
$ws.Range("N74").Value = -3906.8
$ws.Range("H77").Value = 1650.9242
$ws.Range("I77").Value = 1560.2322
$ws.Range("J77").Value = 2158.8
$ws.Range("K77").Value = 7801.161
$ws.Range("L77").Value = 10794
$ws.Range("M77").Value = -3433.161
$ws.Range("N77").Value = -19530
$ws.Range("H132").Value = 3151.1667
$ws.Range("I132").Value = 1660.3462
$ws.Range("K132").Value = 4981.0386
$ws.Range("M132").Value = -2451.0386
$ws.Range("H136").Value = 10446.077
$ws.Range("I136").Value = 12689.9
$ws.Range("J136").Value = 2966.6667
$ws.Range("K136").Value = 38069.7
$ws.Range("L136").Value = 8900.000100000001
$ws.Range("M136").Value = -35519.7
$ws.Range("N136").Value = -14000.0001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 23000
$ws.Range("J4").Value = 10000
$ws.Range("L4").Value = 10000
$ws.Range("N4").Value = -10224
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 623.7646999999999
$ws.Range("I2").Value = 849.25
$ws.Range("J2").Value = 82.59999999999999
$ws.Range("K2").Value = 5095.5
$ws.Range("L2").Value = 495.6
$ws.Range("M2").Value = -4982.5
$ws.Range("N2").Value = -721.5999999999999
$ws.Range("H37").Value = 40700
$ws.Range("J37").Value = 40700
$ws.Range("L37").Value = 122100
$ws.Range("N37").Value = -122324
$ws.Range("H106").Value = 3393.9092
$ws.Range("J106").Value = 3393.9092
$ws.Range("L106").Value = 10181.7276
$ws.Range("N106").Value = -12073.7276
$ws.Range("H138").Value = 13815.1
$ws.Range("I138").Value = 14983.444
$ws.Range("J138").Value = 3300
$ws.Range("K138").Value = 44950.33199999999
$ws.Range("L138").Value = 9900
$ws.Range("M138").Value = -39810.33199999999
$ws.Range("N138").Value = -20180
$ws.Range("H139").Value = 4929.7
$ws.Range("I139").Value = 9824
$ws.Range("J139").Value = 2832.1428
$ws.Range("K139").Value = 29472
$ws.Range("L139").Value = 8496.428400000001
$ws.Range("M139").Value = -24332
$ws.Range("N139").Value = -18776.4284
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 500
$ws.Range("I13").Value = 500
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 500
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -361
$ws.Range("N13").ClearContents()
$ws.Range("H132").Value = 3267.12
$ws.Range("J132").Value = 3748.682
$ws.Range("L132").Value = 11246.046
$ws.Range("N132").Value = -16306.046
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 5571.4287
$ws.Range("J2").Value = 5692.3076
$ws.Range("L2").Value = 5692.3076
$ws.Range("N2").Value = -5916.3076
$ws.Range("H21").Value = 40007
$ws.Range("J21").Value = 40007
$ws.Range("L21").Value = 40007
$ws.Range("N21").Value = -40355
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 39980
$ws.Range("J75").Value = 39980
$ws.Range("L75").Value = 39980
$ws.Range("N75").Value = -41852
$ws.Range("H78").Value = 39980
$ws.Range("J78").Value = 39980
$ws.Range("L78").Value = 119940
$ws.Range("N78").Value = -129300
